$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.215.37'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '1.662.42'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''217.36'
$ws.Range('D6').Value = '''0.5222'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').Value = '''1.003'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '''0.2642'
$ws.Range('E8').Value = '  -0.85%  '
$ws.Range('D9').Value = '''0.06273'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').Value = '''20.78'
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('D11').Value = '''0.07771'
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').Value = '1.760.81'
$ws.Range('E12').Value = '  +5.25%  '
$ws.Range('D13').Value = '''4.469'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = '1.891.03'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '''0.5446'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '0.0₅8148'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '''64.95'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').Value = '26.215.00'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '''4.602'
$ws.Range('E20').Value = '  -2.69%  '
$ws.Range('D21').Value = '''191.39'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').Value = '''10.01'
$ws.Range('E22').Value = '  -2.78%  '
$ws.Range('D23').Value = '''6.003'
$ws.Range('E23').Value = '  -4.04%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').Value = '''138.78'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '''0.1233'
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('D27').Value = '''7.253'
$ws.Range('E27').Value = '  -1.64%  '
$ws.Range('D28').Value = '''16.13'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').Value = '''1.413'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '''0.05947'
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('D31').Value = '''1.273'
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('D32').Value = '''3.530'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('D33').Value = '''3.264'
$ws.Range('E33').Value = '  -3.60%  '
$ws.Range('D34').Value = '''1.577'
$ws.Range('E34').Value = '  -5.97%  '
$ws.Range('D35').Value = '''0.9583'
$ws.Range('E35').Value = '  -4.19%  '
$ws.Range('D36').Value = '''2.420'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = '''2.771'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = '''0.5661'
$ws.Range('E38').Value = '  -6.22%  '
$ws.Range('D39').Value = '''0.01598'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').Value = '''5.959'
$ws.Range('E40').Value = '  -1.15%  '
$ws.Range('D41').Value = '''0.8518'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').Value = '''1.003'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '''100.56'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '1.005.16'
$ws.Range('E44').Value = '  -7.49%  '
$ws.Range('D45').Value = '1.805.51'
$ws.Range('E45').Value = '  -0.37%  '
$ws.Range('D46').Value = '''56.64'
$ws.Range('E46').Value = '  -2.10%  '
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.4346'
$ws.Range('E49').Value = '  +2.70%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''8.000'
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('D51').Value = '''0.05148'
